# Update the date heading and the 25 "three-digit ÷ one-digit" answer
# cells in the worksheet table to the new day's generated values.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-25 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-26 Wednesday", 2)
$d.Content.Find.Execute("228÷3=76, 0", $true, $false, $false, $false, $false, $true, 1, $false, "382÷2=191, 0", 2)
$d.Content.Find.Execute("400÷4=100, 0", $true, $false, $false, $false, $false, $true, 1, $false, "239÷9=26, 5", 2)
$d.Content.Find.Execute("489÷8=61, 1", $true, $false, $false, $false, $false, $true, 1, $false, "447÷7=63, 6", 2)
$d.Content.Find.Execute("418÷6=69, 4", $true, $false, $false, $false, $false, $true, 1, $false, "787÷4=196, 3", 2)
$d.Content.Find.Execute("715÷4=178, 3", $true, $false, $false, $false, $false, $true, 1, $false, "960÷9=106, 6", 2)
$d.Content.Find.Execute("177÷9=19, 6", $true, $false, $false, $false, $false, $true, 1, $false, "552÷6=92, 0", 2)
$d.Content.Find.Execute("344÷7=49, 1", $true, $false, $false, $false, $false, $true, 1, $false, "662÷9=73, 5", 2)
$d.Content.Find.Execute("485÷8=60, 5", $true, $false, $false, $false, $false, $true, 1, $false, "161÷7=23, 0", 2)
$d.Content.Find.Execute("225÷3=75, 0", $true, $false, $false, $false, $false, $true, 1, $false, "443÷9=49, 2", 2)
$d.Content.Find.Execute("842÷3=280, 2", $true, $false, $false, $false, $false, $true, 1, $false, "361÷5=72, 1", 2)
$d.Content.Find.Execute("979÷8=122, 3", $true, $false, $false, $false, $false, $true, 1, $false, "372÷5=74, 2", 2)
$d.Content.Find.Execute("905÷5=181, 0", $true, $false, $false, $false, $false, $true, 1, $false, "643÷5=128, 3", 2)
$d.Content.Find.Execute("852÷4=213, 0", $true, $false, $false, $false, $false, $true, 1, $false, "452÷5=90, 2", 2)
$d.Content.Find.Execute("594÷3=198, 0", $true, $false, $false, $false, $false, $true, 1, $false, "877÷8=109, 5", 2)
$d.Content.Find.Execute("211÷9=23, 4", $true, $false, $false, $false, $false, $true, 1, $false, "204÷4=51, 0", 2)
$d.Content.Find.Execute("598÷9=66, 4", $true, $false, $false, $false, $false, $true, 1, $false, "910÷8=113, 6", 2)
$d.Content.Find.Execute("373÷2=186, 1", $true, $false, $false, $false, $false, $true, 1, $false, "910÷2=455, 0", 2)
$d.Content.Find.Execute("972÷7=138, 6", $true, $false, $false, $false, $false, $true, 1, $false, "718÷6=119, 4", 2)
$d.Content.Find.Execute("434÷6=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "404÷6=67, 2", 2)
$d.Content.Find.Execute("631÷3=210, 1", $true, $false, $false, $false, $false, $true, 1, $false, "871÷8=108, 7", 2)
$d.Content.Find.Execute("489÷4=122, 1", $true, $false, $false, $false, $false, $true, 1, $false, "835÷4=208, 3", 2)
$d.Content.Find.Execute("807÷4=201, 3", $true, $false, $false, $false, $false, $true, 1, $false, "193÷3=64, 1", 2)
$d.Content.Find.Execute("982÷3=327, 1", $true, $false, $false, $false, $false, $true, 1, $false, "491÷2=245, 1", 2)
$d.Content.Find.Execute("134÷5=26, 4", $true, $false, $false, $false, $false, $true, 1, $false, "348÷5=69, 3", 2)
$d.Content.Find.Execute("373÷9=41, 4", $true, $false, $false, $false, $false, $true, 1, $false, "915÷9=101, 6", 2)
